$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 74
$ws.Range("I2").Value = 170
$ws.Range("J2").Value = 705
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 178
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 99
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 8
$ws.Range("S2").Value = 74
$ws.Range("T2").Value = 119
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 1105
$ws.Range("X2").Value = 1129
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 13
